$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time Tracker(Time Table)")

# Copy formats of the last existing data row (89) down to the 4 new rows (90-93)
$ws.Range("A89:W89").Copy()
$ws.Range("A90:W93").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights (matches ht="19" in the target rows)
$ws.Rows.Item(90).RowHeight = 19
$ws.Rows.Item(91).RowHeight = 19
$ws.Rows.Item(92).RowHeight = 19
$ws.Rows.Item(93).RowHeight = 19

# Row 90 - wash / 1hour
$ws.Range("A90").Value2 = 44382
$ws.Range("B90").Value2 = 0.45833333333333331
$ws.Range("C90").Value2 = 44382
$ws.Range("D90").Value2 = 0.5
$ws.Range("E90").Value2 = "root lim"
$ws.Range("F90").Value2 = "RootLimSecretary"
$ws.Range("G90").Value2 = "wash"
$ws.Range("H90").Value2 = "1hour"
$ws.Range("I90").Value2 = "필수불가결한일(생활에)"
$ws.Range("J90").Value2 = "wash"
$ws.Range("K90").Value2 = "1hour"

# Row 91 - research useRef / 2.5 hour
$ws.Range("A91").Value2 = 44382
$ws.Range("B91").Value2 = 0.52083333333333337
$ws.Range("C91").Value2 = 44382
$ws.Range("D91").Value2 = 0.60416666666666663
$ws.Range("E91").Value2 = "root lim"
$ws.Range("F91").Value2 = "ttkmw"
$ws.Range("G91").Value2 = "research useRef"
$ws.Range("H91").Value2 = "2.5 hour"

# Row 92 - eat, stretch / 1 hour
$ws.Range("A92").Value2 = 44382
$ws.Range("B92").Value2 = 0.60416666666666663
$ws.Range("C92").Value2 = 44382
$ws.Range("D92").Value2 = 0.625
$ws.Range("E92").Value2 = "root lim"
$ws.Range("F92").Value2 = "RootLimSecretary"
$ws.Range("G92").Value2 = "eat, stretch"
$ws.Range("H92").Value2 = "1 hour"
$ws.Range("I92").Value2 = "필수불가결한일(생활에)"
$ws.Range("J92").Value2 = "eat"

# Row 93 - develop table / 3 hour (End Time left blank, general format not time)
$ws.Range("A93").Value2 = 44382
$ws.Range("B93").Value2 = 0.62569444444444444
$ws.Range("C93").Value2 = 44382
$ws.Range("L90").Copy()
$ws.Range("D93").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D93").ClearContents()
$ws.Range("E93").Value2 = "root lim"
$ws.Range("F93").Value2 = "ttkmw"
$ws.Range("G93").Value2 = "develop table"
$ws.Range("H93").Value2 = "3 hour"

# The table ("표3") auto-expands to cover the newly added rows
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A5:W93"))

Write-Output "done"
